$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.428.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.979.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.80%  "
$ws.Range("E4").Value = "  -0.91%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "502.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.35%  "
$ws.Range("E8").Value = "  +8.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.32"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +11.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.105"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.351"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.22%  "
$ws.Range("E12").Value = "  +2.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.491.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +12.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.427.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000150"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +15.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.978.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.57%  "
$ws.Range("E18").Value = "  +11.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.55%  "
$ws.Range("E20").Value = "  +11.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.470"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.162"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0890"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +13.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +12.01%  "
$ws.Range("E30").Value = "  +12.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.90%  "
$ws.Range("E32").Value = "  +9.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "156.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.18%  "
$ws.Range("E35").Value = "  +8.29%  "
$ws.Range("E36").Value = "  +5.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.53"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.40%  "
$ws.Range("E38").Value = "  +11.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "23.05"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.014.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "36.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.640"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.241.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.989"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.58%  "
$ws.Range("E46").Value = "  +7.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +26.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0235"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +11.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.77%  "
$ws.Range("E51").Value = "  +8.89%  "

Write-Host "Updated cryptos list"
